# Template import YCDC: remove the "Chi tiết / Đơn vị (*)" (unit) column.
# The column D header ("Chi tiết / Đơn vị (*)") and its "Cái" values are no
# longer needed, so the whole column is deleted and everything to its right
# (Từ kho / Đến kho / Từ LSX / Đến LSX) shifts one slot to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns(4).Delete()

# Restore the selection the author left the sheet on.
$ws.Range("G12").Select() | Out-Null
